$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width: 16.42578125 -> 15.42578125 (stored OOXML width units)
# ColumnWidth setter applies a +0.8333333333333334 offset internally and
# snaps to 1/6-character increments, so the nearest achievable stored
# width to 15.42578125 is 15.5 (input 14.666666666666666).
$ws.Columns(1).ColumnWidth = 14.666666666666666

# Update cell values in column A
$ws.Range("A1").Value = 0.32703571387048669
$ws.Range("A2").Value = -0.032854260718867323
$ws.Range("A3").Value = -0.0039999999677728937
$ws.Range("A4").Value = -0.0079999999426103585
$ws.Range("A5").Value = -0.0029999999637277952
$ws.Range("A6").Value = -0.001999999955168974
$ws.Range("A7").Value = -0.0099999999131421546
$ws.Range("A8").Value = -0.0099999999113058458
$ws.Range("A9").Value = -0.0019999999501836285
$ws.Range("A10").Value = -0.0019999999482536168
$ws.Range("A11").Value = -0.0029999999429959345
$ws.Range("A12").Value = -0.0034999999400495696
$ws.Range("A13").Value = -0.0034999999385822989
$ws.Range("A14").Value = -0.007999999915139
$ws.Range("A15").Value = 0.010250906017729378
$ws.Range("A16").Value = -0.0019999999454625161
$ws.Range("A17").Value = -0.0019999999446200789
$ws.Range("A18").Value = -0.0039999999341215897
$ws.Range("A19").Value = -0.0039999999764126493
$ws.Range("A20").Value = -0.018911512655085616
$ws.Range("A21").Value = -0.060300293680654704
$ws.Range("A22").Value = -0.0039999999622608584
$ws.Range("A23").Value = -0.004999999961114554
$ws.Range("A24").Value = -0.019999999880558228
$ws.Range("A25").Value = -0.019999999879015462
$ws.Range("A26").Value = -0.0024999999511177151
$ws.Range("A27").Value = -0.002499999948787579
$ws.Range("A28").Value = -0.0019999999403559343
$ws.Range("A29").Value = -0.0069999999069088048
$ws.Range("A30").Value = -0.059999999630838907
$ws.Range("A31").Value = -0.0069999999006817859
$ws.Range("A32").Value = -0.0099999998849042981
$ws.Range("A33").Value = -0.0039999999155817534
